$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.363.89"
$ws.Range("E2").Value = "  +0.00%  "
$ws.Range("D3").Value = "1.843.55"
$ws.Range("E3").Value = "  -0.20%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "'240.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.10%  "
$ws.Range("D6").Value = "'0.6303"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("D7").Value = "'0.9997"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "'0.07458"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.74%  "
$ws.Range("D9").Value = "'0.2909"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("D10").Value = "'24.95"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.14%  "
$ws.Range("E11").Value = "  +0.01%  "
$ws.Range("D12").Value = "1.845.79"
$ws.Range("E12").Value = "  -0.07%  "
$ws.Range("D13").Value = "'4.990"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.58%  "
$ws.Range("D14").Value = "'0.6793"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.13%  "
$ws.Range("D16").Value = "'82.21"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.02%  "
$ws.Range("D17").Value = "'6.282"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.68%  "
$ws.Range("D18").Value = "29.382.40"
$ws.Range("E18").Value = "  -0.01%  "
$ws.Range("D19").Value = "'229.54"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.26%  "
$ws.Range("D20").Value = "'12.33"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.08%  "
$ws.Range("D21").Value = "'0.9999"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("D22").Value = "'7.415"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.45%  "
$ws.Range("E23").Value = "  -0.15%  "
$ws.Range("D24").Value = "'158.42"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.36%  "
$ws.Range("D25").Value = "'8.491"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.89%  "
$ws.Range("E26").Value = "  -2.31%  "
$ws.Range("E27").Value = "  -0.90%  "
$ws.Range("D28").Value = "'0.06539"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +14.92%  "
$ws.Range("D29").Value = "'1.441"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.13%  "
$ws.Range("D30").Value = "'1.487"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.73%  "
$ws.Range("D31").Value = "'4.074"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.43%  "
$ws.Range("D32").Value = "'4.061"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.47%  "
$ws.Range("D33").Value = "'1.841"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.32%  "
$ws.Range("E34").Value = "  -1.10%  "
$ws.Range("D35").Value = "'0.6977"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.43%  "
$ws.Range("D36").Value = "'2.578"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.34%  "
$ws.Range("D37").Value = "'0.01856"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.54%  "
$ws.Range("D38").Value = "'2.817"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.33%  "
$ws.Range("D39").Value = "1.248.64"
$ws.Range("E39").Value = "  -0.13%  "
$ws.Range("D40").Value = "'6.780"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.64%  "
$ws.Range("D41").Value = "'0.9317"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.97%  "
$ws.Range("D42").Value = "'0.9994"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.05%  "
$ws.Range("D43").Value = "1.999.66"
$ws.Range("E43").Value = "  -0.62%  "
$ws.Range("D44").Value = "'100.81"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.59%  "
$ws.Range("D45").Value = "'65.59"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.49%  "
$ws.Range("E46").Value = "  +3.48%  "
$ws.Range("D47").Value = "'7.067"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.21%  "
$ws.Range("D48").Value = "'1.716"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.94%  "
$ws.Range("D49").Value = "'8.994"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.28%  "
$ws.Range("E50").Value = "  -1.46%  "
$ws.Range("D51").Value = "'0.3902"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.32%  "
